$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand used range / header row (row 1): B1, C1, D1
$ws.Range("B1").Value = "Хэш"
$ws.Range("C1").Value = "Время обработки"
$ws.Range("D1").Value = "Хэммингово расстояние"

# Format column B (rows 2-63) as Text so the 68-character binary-looking
# phash strings are stored verbatim and not coerced into scientific-notation numbers.
$ws.Range("B2:B63").NumberFormat = "@"

# Data rows 2-63: hash (B), processing time (C), hamming distance (D)
$data = @(
    @(2, "1011100100000000110001000000000000000000000000000001000000000000", 0.031247, 0),
    @(3, "1011100100000000110000000000000000000000000000000000000000000000", 0.015618, 2),
    @(4, "1011100100000000110001000000000000000000000000000001000000000000", 0.015627, 0),
    @(5, "1011100100000000110001000000000000000000000000000001000000000000", 0.03126, 0),
    @(6, "1011100100000000110000000000000000000000000000000000000000000000", 0, 2),
    @(7, "1011100100000000110000000000000000000000000000000000000000000000", 0, 2),
    @(8, "1011100100000000110000000000000000000000000000000000000000000000", 0, 2),
    @(9, "1011100100000000110000000000000000000000000000000000000000000000", 0, 2),
    @(10, "1011100100000000110000000000000000000000000000000000000000000000", 0, 2),
    @(11, "1011100100000000110000000000000000000000000000000000000000000000", 0, 2),
    @(12, "1011100100000000110001000000000000000000000000000001000000000000", 0, 0),
    @(13, "1011100100000000110001000000000000000000000000000001000000000000", 0, 0),
    @(14, "1011100100000000110001000000000000000000000000000001000000000000", 0, 0),
    @(15, "1011100100000000110001000000000000000000000000000001000000000000", 0.015617, 0),
    @(16, "1011100100000000110001000000000000000000000000000001000000000000", 0, 0),
    @(17, "1011100100000000110001000000000000000000000000000001000000000000", 0, 0),
    @(18, "1011100100000000110001000000000000000000000000000001000000000000", 0, 0),
    @(19, "1011100100000000110001000000000000000000000000000001000000000000", 0, 0),
    @(20, "1011100100000000110001000000000000000000000000000001000000000000", 0, 0),
    @(21, "1011100100000000110001000000000000000000000000000001000000000000", 0, 0),
    @(22, "1011100100000000110001000000000010000000000000000001000000000000", 0, 1),
    @(23, "1011100100000000110001000000000000000000000000000001000000000000", 0.015624, 0),
    @(24, "1011100100000000110001000000000000000000000000000001000000000000", 0.015632, 0),
    @(25, "1011100100000000110001000000000000000000000000000001000000000000", 0.015556, 0),
    @(26, "1011101100000000110000000000000000000000000000000001000000000000", 0.015702, 2),
    @(27, "1011101100000000110000000000000000000000000000000000000000000000", 0.015618, 3),
    @(28, "1011101100000000110000000000000001000000000000000000000000000000", 0.015631, 4),
    @(29, "1011001100000000110000000000000000000001000000000000000000000000", 0, 5),
    @(30, "1011001100000000110000000000000000000001000000000000000000000000", 0, 5),
    @(31, "1011001100000000110000000000000000000001000000000000000000000000", 0.015625, 5),
    @(32, "1011001100000000110000000000000000000001000000000000000000000000", 0.015634, 5),
    @(33, "1011001100000000110000000000000000000001000000000000000000000000", 0.015625, 5),
    @(34, "1010000000000000100000000000000000000000000000000000000000000000", 0, 6),
    @(35, "1110101101000000101001000001000010100100000000001000011000000000", 0.01562, 14),
    @(36, "1110100100000000100000000000000010000010100000000000111000000000", 0, 11),
    @(37, "1110100000000000100100001000000000000000000000000000000001000000", 0, 9),
    @(38, "1110110001100000101000001000100110000000000000001000000000000000", 0, 15),
    @(39, "1011100000000000110000000000000000000000000000000000000000000000", 0, 3),
    @(40, "1010100100000000100000000000000010000010000000001000000000000000", 0, 7),
    @(41, "1011100000000000110000001000100000000000100000000000000000000000", 0.015623, 6),
    @(42, "1011000000000000100000000000000000000000000000000000000000000000", 0, 5),
    @(43, "1011100100000000110000000000000010000000000000000001000000000000", 0, 2),
    @(44, "1011000000110000110000001100000010000000000000000000000000000000", 0.015559, 9),
    @(45, "1011101100000000110000000000000010000000000000000001000000000000", 0, 3),
    @(46, "1011100100000000110001001000000000010000000000000000000000000000", 0, 3),
    @(47, "1011100000000000110000001000000000000000000000000000000000000000", 0.015678, 4),
    @(48, "1011001100000000110000000000000010000000000000000000000000000000", 0, 5),
    @(49, "1111100100000000100001000001000010000010000000001000000000000000", 0.015585, 7),
    @(50, "1110100100000000110011000000000000000000000000000010000000000000", 0, 5),
    @(51, "1011101100000000110000001000000000000000000000000001000000000000", 0.015697, 3),
    @(52, "1011101100000000110000001000000000000000000000000001000000000000", 0.015588, 3),
    @(53, "1011100100000000110000000000000000000000000000000001000000000000", 0, 1),
    @(54, "1010001000000000110000000000000000000000000000000000000000000000", 0.015647, 6),
    @(55, "1010000000000000100000000000000000010000000000000000000000000000", 0, 7),
    @(56, "1100111000010000101100000000000000000000100000000000000000000000", 0.031191, 13),
    @(57, "1011000000000000100000100000000000000000000000001000000000000000", 0, 7),
    @(58, "1011101100000000110001000000000000010010000000000001000000000000", 0.03124, 3),
    @(59, "1010000000000000100000100000000010001000000000001000000000000000", 0.015624, 10),
    @(60, "1010100000010000110000001100000000000010100000011000000000010000", 0.046875, 12),
    @(61, "1011101001000000111000000000000000000000100000000000010100000000", 0.046883, 9),
    @(62, "1011101010000000110000000000000011000000000000000000001000010000", 0.031259, 9),
    @(63, "1011000000000000100000000000000000000000000000000000000000000000", 0.062483, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
